# Retraining the model for Astro
# Refresh the rolling Production forecast table (A2:D170) on Sheet1 with
# the newly retrained model's output: date (A), hour-of-day (B),
# predicted value (C) and the Date+Hour lookup label (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: serial date for each forecast row
$dates  = @(46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46046,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46047,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46048,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46049,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46050,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46051,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46052,46053,46053,46053,46053,46053,46053,46053,46053,46053)
# Column B: hour of day (1-24) for each forecast row
$hours  = @(9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9)
# Column C: predicted value for each forecast row
$vals   = @(0.068,0.068,0.354,0.5570000000000001,0.45,0.386,0.373,0.264,0.13,0.032,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.018,0.178,0.356,0.447,0.449,0.404,0.327,0.179,0.031,0,0,0.012,0.012,0.012,0.012,0,0,0,0,0,0,0,0,0,0,0.11,0.249,0.303,0.31,0.27,0.186,0.115,0.041,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.014,0.13,0.301,0.474,0.572,0.5669999999999999,0.496,0.321,0.117,0,0,0,0,0,0.012,0,0,0,0,0,0,0,0,0,0.08599999999999999,0.37,0.786,0.854,0.833,0.859,0.629,0.412,0.147,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.078,0.365,0.589,0.762,0.835,0.888,0.666,0.424,0.135,0,0,0,0,0.012,0.222,0.022,0.026,0.016,0.016,0.016,0.016,0.017,0.017,0.083,0.192,0.622,1.049,1.434,0.975,0.986,0.901,0.476,0.16,0.017,0.193,0.212,0.195,0.2,0.2,0.022,0.027,0.023,0.023,0.02,0.02,0.02,0.02,0.101)
# Column D: "Lookup" label combining date and hour
$labels = @("24.01.20269","24.01.202610","24.01.202611","24.01.202612","24.01.202613","24.01.202614","24.01.202615","24.01.202616","24.01.202617","24.01.202618","24.01.202619","24.01.202620","24.01.202621","24.01.202622","24.01.202623","24.01.202624","25.01.20261","25.01.20262","25.01.20263","25.01.20264","25.01.20265","25.01.20266","25.01.20267","25.01.20268","25.01.20269","25.01.202610","25.01.202611","25.01.202612","25.01.202613","25.01.202614","25.01.202615","25.01.202616","25.01.202617","25.01.202618","25.01.202619","25.01.202620","25.01.202621","25.01.202622","25.01.202623","25.01.202624","26.01.20261","26.01.20262","26.01.20263","26.01.20264","26.01.20265","26.01.20266","26.01.20267","26.01.20268","26.01.20269","26.01.202610","26.01.202611","26.01.202612","26.01.202613","26.01.202614","26.01.202615","26.01.202616","26.01.202617","26.01.202618","26.01.202619","26.01.202620","26.01.202621","26.01.202622","26.01.202623","26.01.202624","27.01.20261","27.01.20262","27.01.20263","27.01.20264","27.01.20265","27.01.20266","27.01.20267","27.01.20268","27.01.20269","27.01.202610","27.01.202611","27.01.202612","27.01.202613","27.01.202614","27.01.202615","27.01.202616","27.01.202617","27.01.202618","27.01.202619","27.01.202620","27.01.202621","27.01.202622","27.01.202623","27.01.202624","28.01.20261","28.01.20262","28.01.20263","28.01.20264","28.01.20265","28.01.20266","28.01.20267","28.01.20268","28.01.20269","28.01.202610","28.01.202611","28.01.202612","28.01.202613","28.01.202614","28.01.202615","28.01.202616","28.01.202617","28.01.202618","28.01.202619","28.01.202620","28.01.202621","28.01.202622","28.01.202623","28.01.202624","29.01.20261","29.01.20262","29.01.20263","29.01.20264","29.01.20265","29.01.20266","29.01.20267","29.01.20268","29.01.20269","29.01.202610","29.01.202611","29.01.202612","29.01.202613","29.01.202614","29.01.202615","29.01.202616","29.01.202617","29.01.202618","29.01.202619","29.01.202620","29.01.202621","29.01.202622","29.01.202623","29.01.202624","30.01.20261","30.01.20262","30.01.20263","30.01.20264","30.01.20265","30.01.20266","30.01.20267","30.01.20268","30.01.20269","30.01.202610","30.01.202611","30.01.202612","30.01.202613","30.01.202614","30.01.202615","30.01.202616","30.01.202617","30.01.202618","30.01.202619","30.01.202620","30.01.202621","30.01.202622","30.01.202623","30.01.202624","31.01.20261","31.01.20262","31.01.20263","31.01.20264","31.01.20265","31.01.20266","31.01.20267","31.01.20268","31.01.20269")

$firstRow = 2
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $hours[$i]
    $ws.Cells.Item($row, 3).Value = $vals[$i]
    $ws.Cells.Item($row, 4).Value = $labels[$i]
}

